# Budget_Dashboard / Settings sheet: add a "Starting Year" setting row
# inside the existing "General" box, growing the box by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# ---------------------------------------------------------------------
# 1) Make room: the box currently occupies rows 5-8 (header + 2 empty
#    rows + bottom border row). The target has it occupying rows 6-10
#    (header + empty row + new "Starting Year" row + 2 more rows).
#    Insert a row above the header (shifts everything down by 1) and a
#    second row before the final (bottom-border) row.
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(9).Insert()

# xlEdge constants
$xlLeft = 7
$xlTop = 8
$xlBottom = 9
$xlRight = 10
$xlContinuous = 1
$xlNone = -4142
$xlThin = 2
$xlCenter = -4108

function Clear-Borders($rng) {
    $rng.Borders.Item($xlLeft).LineStyle = $xlNone
    $rng.Borders.Item($xlTop).LineStyle = $xlNone
    $rng.Borders.Item($xlBottom).LineStyle = $xlNone
    $rng.Borders.Item($xlRight).LineStyle = $xlNone
}

function Set-Borders($rng, $left, $top, $right, $bottom) {
    Clear-Borders $rng
    if ($left) {
        $rng.Borders.Item($xlLeft).LineStyle = $xlContinuous
        $rng.Borders.Item($xlLeft).Weight = $xlThin
    }
    if ($top) {
        $rng.Borders.Item($xlTop).LineStyle = $xlContinuous
        $rng.Borders.Item($xlTop).Weight = $xlThin
    }
    if ($right) {
        $rng.Borders.Item($xlRight).LineStyle = $xlContinuous
        $rng.Borders.Item($xlRight).Weight = $xlThin
    }
    if ($bottom) {
        $rng.Borders.Item($xlBottom).LineStyle = $xlContinuous
        $rng.Borders.Item($xlBottom).Weight = $xlThin
    }
}

# fontThemeColor / fillThemeColor use the COM XlThemeColor numbering,
# which (empirically, via this host) maps to the OOXML <color theme="N"/>
# index as: com 1->1, 2->0, 3->3, 4->2, 5->4, 6->5, ...
function Set-CellFormat($rng, $bold, $italic, $fontThemeColor, $fillThemeColor, $hcenter) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Bold = $bold
    $rng.Font.Italic = $italic
    if ($fontThemeColor -ne $null) {
        $rng.Font.ThemeColor = $fontThemeColor
    }
    if ($fillThemeColor -ne $null) {
        $rng.Interior.Pattern = 1
        $rng.Interior.ThemeColor = $fillThemeColor
    }
    if ($hcenter) {
        $rng.HorizontalAlignment = $xlCenter
    }
}

# ---------------------------------------------------------------------
# 2) Row 6: "General" header bar (unchanged look, now one row lower)
# ---------------------------------------------------------------------
$rngHeader = $ws.Range("C6:H6")
Set-CellFormat $rngHeader $false $false 2 5 $true
Clear-Borders $rngHeader
$ws.Range("C6").Value2 = "General"

# ---------------------------------------------------------------------
# 3) Row 7: empty spacer row directly under the header (box side walls)
# ---------------------------------------------------------------------
Set-CellFormat $ws.Range("C7") $false $false 2 5 $true
Set-Borders $ws.Range("C7") $true $false $false $true

Set-CellFormat $ws.Range("E7:G7") $false $false 2 5 $true
Set-Borders $ws.Range("E7:G7") $false $false $false $true

Set-CellFormat $ws.Range("H7") $false $false 2 5 $true
Set-Borders $ws.Range("H7") $false $false $true $true

$ws.Range("D7").ClearContents()

# ---------------------------------------------------------------------
# 4) Row 8: the new "Starting Year" content row
# ---------------------------------------------------------------------
Set-CellFormat $ws.Range("C8") $false $false 2 5 $true
Set-Borders $ws.Range("C8") $true $true $false $false

Set-CellFormat $ws.Range("D8") $true $false $null 4 $true
Set-Borders $ws.Range("D8") $true $false $false $false
$ws.Range("D8").Value2 = "Starting Year:"

Set-CellFormat $ws.Range("E8") $false $false $null 2 $true
Set-Borders $ws.Range("E8") $false $false $false $false
$ws.Range("E8").Value2 = 2023

Set-CellFormat $ws.Range("F8:G8") $false $true 4 4 $true
Set-Borders $ws.Range("F8:G8") $false $false $false $false
$ws.Range("F8").Value2 = "Set the starting year (yyyy) once at the beginning and do not chnge again"
$ws.Range("F8:G8").Merge()

Set-CellFormat $ws.Range("H8") $false $false $null 4 $true
Set-Borders $ws.Range("H8") $false $false $true $false

# ---------------------------------------------------------------------
# 5) Row 9: empty spacer row under the content row
# ---------------------------------------------------------------------
Set-CellFormat $ws.Range("C9") $false $false 2 5 $true
Set-Borders $ws.Range("C9") $true $true $true $true

Set-CellFormat $ws.Range("D9") $false $false $null 4 $true
Set-Borders $ws.Range("D9") $true $false $false $false

Set-CellFormat $ws.Range("E9:G9") $false $false $null 4 $true
Set-Borders $ws.Range("E9:G9") $false $false $false $false

Set-CellFormat $ws.Range("H9") $false $false $null 4 $true
Set-Borders $ws.Range("H9") $false $false $true $false

# ---------------------------------------------------------------------
# 6) Row 10: bottom border row closing the box
# ---------------------------------------------------------------------
Set-CellFormat $ws.Range("C10") $false $false 2 5 $true
Set-Borders $ws.Range("C10") $true $true $true $true

Set-CellFormat $ws.Range("D10") $false $false $null 4 $true
Set-Borders $ws.Range("D10") $true $false $false $true

Set-CellFormat $ws.Range("E10:G10") $false $false $null 4 $true
Set-Borders $ws.Range("E10:G10") $false $false $false $true

Set-CellFormat $ws.Range("H10") $false $false $null 4 $true
Set-Borders $ws.Range("H10") $false $false $true $true

# ---------------------------------------------------------------------
# 7) Column widths: split column G off with a wider width for the
#    help text, leave D:F at their existing width.
# ---------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 50

# ---------------------------------------------------------------------
# 8) Selection / view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------
$ws.Range("G9").Select()

Write-Host "done"
